# Generate Report for Handback
# Adds a new handback entry (87f810d5-e8df-4903-b300-3720d7e05e96) as row 4
# to the Overview, zh-cn and de-de tables.

$wb = $excel.ActiveWorkbook

$fileGuid = "87f810d5-e8df-4903-b300-3720d7e05e96"
$mdName = "$fileGuid.md"
$mdDisplayOverview = "e2e\$fileGuid.md"

# ---------------------------------------------------------------------------
# Sheet "Overview": append summary row
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A4").Value = $mdName
$wsOverview.Range("B4").Value = $mdDisplayOverview
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("G4").Value = "2016-08-30 18:55:52"

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6f7a7c9d9a0a6e6a9f3b6a1c9e8f7d6c5b4a3928/e2e/$mdName", [Type]::Missing, [Type]::Missing, $mdDisplayOverview) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn": append detail row
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$zhXlf = "$fileGuid.53fb24ab51358fb434bf3ed4ff613733c4f2bd7c.zh-cn.xlf"

$wsZhCn.Range("A4").Value = $mdName
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "True"
$wsZhCn.Range("G4").Value = $zhXlf
$wsZhCn.Range("H4").Value = "2016-08-30 18:55:47"
$wsZhCn.Range("I4").Value = $mdName
$wsZhCn.Range("J4").Value = $zhXlf
$wsZhCn.Range("K4").Value = "2016-08-30 18:56:16"
$wsZhCn.Range("M4").Value = "True"
$wsZhCn.Range("O4").Value = "False"

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6f7a7c9d9a0a6e6a9f3b6a1c9e8f7d6c5b4a3928/e2e/$mdName", [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/1a2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b/e2e/$mdName", [Type]::Missing, [Type]::Missing, $mdName) | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de": append detail row
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$deXlf = "$fileGuid.53fb24ab51358fb434bf3ed4ff613733c4f2bd7c.de-de.xlf"

$wsDeDe.Range("A4").Value = $mdName
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "True"
$wsDeDe.Range("G4").Value = $deXlf
$wsDeDe.Range("H4").Value = "2016-08-30 18:55:52"
$wsDeDe.Range("I4").Value = $mdName
$wsDeDe.Range("J4").Value = $deXlf
$wsDeDe.Range("K4").Value = "2016-08-30 18:56:23"
$wsDeDe.Range("M4").Value = "True"
$wsDeDe.Range("O4").Value = "False"

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6f7a7c9d9a0a6e6a9f3b6a1c9e8f7d6c5b4a3928/e2e/$mdName", [Type]::Missing, [Type]::Missing, $mdName) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b1c/e2e/$mdName", [Type]::Missing, [Type]::Missing, $mdName) | Out-Null

Write-Host "Handback report row added for $fileGuid"
